$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (fuzzy-matched crosswalk row)
$ws.Range("A2").Value = "AN22-000224"
$ws.Range("B2").Value = "Clean Up"
$ws.Range("D2").Value = "Chrissy Bruno"
$ws.Range("F2").Value = "3301 Old York Road"

# Date-like text columns: force Text number format before assignment so the
# values stay literal strings (matching source data) instead of being
# auto-parsed into date serials, then reset the style back to Normal so no
# stray number-format style sticks to the cell.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "02-25-2022"
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "02-25-2022"
$ws.Range("H2").Style = "Normal"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "03-31-2022"
$ws.Range("L2").Style = "Normal"

$ws.Range("O2").Value = "Minor Removal Project"
$ws.Range("Q2").Value = "Mary McLeod Bethune School"
$ws.Range("R2").Value = "School District of Philadelphia"
$ws.Range("S2").Value = "440 North Broad Street Philadelphia, PA"
$ws.Range("T2").Value = "PEPPER ENVIRONMENTAL SERVICES"
$ws.Range("AB2").Value = 35

# Remove rows 3 and 4 entirely (data no longer present after fuzzy match)
$ws.Rows("3:4").Delete()
